$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 12.0407021
$ws.Range("D3").Value = 49.681352
